$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.350.02"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "3.665.88"
$ws.Range("E3").Value = "  -3.27%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'596.26"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "'165.74"
$ws.Range("E6").Value = "  -4.23%  "
$ws.Range("D7").Value = "3.664.11"
$ws.Range("E7").Value = "  -3.29%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.531"
$ws.Range("E9").Value = "  +1.11%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("D11").Value = "'6.28"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("D13").Value = "'37.79"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "4.280.04"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("D16").Value = "3.664.46"
$ws.Range("E16").Value = "  -3.38%  "
$ws.Range("D17").Value = "68.182.39"
$ws.Range("E17").Value = "  +0.24%  "
$ws.Range("D18").Value = "'7.19"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("E20").Value = "  +6.05%  "
$ws.Range("D21").Value = "'489.32"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("E23").Value = "  -1.81%  "
$ws.Range("D24").Value = "'84.29"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("E26").Value = "  -4.35%  "
$ws.Range("E27").Value = "  -0.62%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'10.03"
$ws.Range("E29").Value = "  -1.79%  "
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("D31").Value = "'7.81"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("E32").Value = "  -2.82%  "
$ws.Range("D33").Value = "'31.14"
$ws.Range("E33").Value = "  -4.48%  "
$ws.Range("D34").Value = "3.806.43"
$ws.Range("E34").Value = "  -3.22%  "
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").Value = "3.607.56"
$ws.Range("E36").Value = "  -3.26%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("D38").Value = "'0.992"
$ws.Range("E38").Value = "  -1.29%  "
$ws.Range("E39").Value = "  -1.04%  "
$ws.Range("D40").Value = "'0.131"
$ws.Range("E40").Value = "  -3.94%  "
$ws.Range("E41").Value = "  -1.89%  "
$ws.Range("D42").Value = "'48.86"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").Value = "'428.68"
$ws.Range("E43").Value = "  -5.26%  "
$ws.Range("D44").Value = "'1.95"
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("E45").Value = "  -3.31%  "
$ws.Range("E46").Value = "  +1.05%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'40.09"
$ws.Range("E48").Value = "  -3.23%  "
$ws.Range("D49").Value = "'141.42"
$ws.Range("E49").Value = "  +1.06%  "
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("D51").Value = "2.716.48"
$ws.Range("E51").Value = "  -3.80%  "
